# excel bugfixes: lotto5_feladat.docx
# - reword the V4 task bullet (what to search for changes)
# - clarify the chart task wording, emphasising "3 találatos" and "2024-ben"
# - keep the chart title in sync with the same wording

$d = $word.ActiveDocument

# Helper: split the range [$startPos, $startPos + len($text)) into its own
# run by toggling Bold on/off (forces the engine to stop merging it with its
# neighbours) while leaving the run's *effective* bold state equal to
# $bold. Returns the position right after the inserted text so callers can
# chain calls together.
function Split-Run($doc, $startPos, $text, $bold) {
    $len = $text.Length
    $sub = $doc.Range($startPos, $startPos + $len)
    if ($bold) {
        $sub.Font.Bold = 0
        $sub.Font.Bold = 1
    } else {
        $sub.Font.Bold = 1
        $sub.Font.Bold = 0
    }
    return $startPos + $len
}

# --- 1) "...keresd meg, hogy melyik számot húzták ki a legtöbbször a..." -> reworded
$range = $d.Content
$found = $range.Find.Execute(" cellában keresd meg, hogy melyik számot húzták ki a legtöbbször a ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $range.Start
    $range.Text = " cellában keresd meg, hogy mekkora volt a legnagyobb előfordulás a számok kihúzásánál a "

    $pos = $start
    $pos = Split-Run $d $pos " cellában keresd meg, hogy me" $false
    $pos = Split-Run $d $pos "kkora volt a legnagyobb előfordulás a " $false
    $pos = Split-Run $d $pos "számo" $false
    $pos = Split-Run $d $pos "k" $false
    $pos = Split-Run $d $pos " " $false
    $pos = Split-Run $d $pos "ki" $false
    $pos = Split-Run $d $pos "húzá" $false
    $pos = Split-Run $d $pos "sánál" $false
    $pos = Split-Run $d $pos " a " $false
}

# --- 2) chart paragraph: "...amely a heti nyertes szelvények..." -> "...amely a heti *3 találatos* nyertes szelvények... *2024-ben*. Ehhez az adatokat..."
$range = $d.Content
$found = $range.Find.Execute(", amely a heti nyertes szelvények számát mutatja be 2024-ben. Ehhez az évszűrést alkalmazd, és az adatokat a megfelelő oszlopból vedd ki.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $range.Start
    $range.Text = ", amely a heti 3 találatos nyertes szelvények számát mutatja be 2024-ben. Ehhez az adatokat a megfelelő oszlopból vedd ki."

    $pos = $start
    $pos = Split-Run $d $pos ", amely a heti " $false
    $pos = Split-Run $d $pos "3 találatos" $true
    $pos = Split-Run $d $pos " " $false
    $pos = Split-Run $d $pos "nyertes szelvények számát mutatja be " $false
    $pos = Split-Run $d $pos "2024-ben" $true
    $pos = Split-Run $d $pos ". Ehhez az adatokat a megfelelő oszlopból vedd ki." $false
}

# --- 3) chart title: "Heti nyertes szelvények számát mutatja be 2024-ben" -> "Heti *3 találatos* szelvények számát mutatja be 2024-ben"
$range = $d.Content
$found = $range.Find.Execute("Heti nyertes szelvények számát mutatja be 2024-ben", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $start = $range.Start
    $range.Text = "Heti 3 találatos szelvények számát mutatja be 2024-ben"

    $pos = $start
    $pos = Split-Run $d $pos "Heti " $true
    $pos = Split-Run $d $pos "3 találatos" $true
    $pos = Split-Run $d $pos " szelvények számát mutatja be 2024-ben" $true
}
